$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatisticsRoman")

$ws.Range("B11").Value = 6035
$ws.Range("B12").Value = 1.44
$ws.Range("B20").Value = 1.07

$ws.Activate() | Out-Null
$ws.Range("G8").Select() | Out-Null
